$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G3").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G4").Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G5").Value = 'hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G6").Value = 'mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Range("G7").Value = 'hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G8").Value = 'backup@backdoor.com, manar.montaser@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G9").Value = 'hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G10").Value = 'drsarahmw@gmail.com, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G11").Value = 'hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G13").Value = 'Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G14").Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G15").Value = 'aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G16").Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G17").Value = 'dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G19").Value = 'Sarah.Mahdy@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Range("G20").Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G22").Value = 'nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Range("G24").Value = 'Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G25").Value = 'Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G27").Value = 'Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G28").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range("G29").Value = 'neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G30").Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G31").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G32").Value = 'hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G33").Value = 'hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G34").Value = 'mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Range("G35").Value = 'hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G36").Value = 'backup@backdoor.com, manar.montaser@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G37").Value = 'hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G38").Value = 'drsarahmw@gmail.com, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G39").Value = 'hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G41").Value = 'Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G42").Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G43").Value = 'aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G44").Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G45").Value = 'dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G47").Value = 'Sarah.Mahdy@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Range("G48").Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G50").Value = 'nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Range("G52").Value = 'Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G53").Value = 'Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G55").Value = 'Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G56").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range("G57").Value = 'neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
